$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")
$ws.Activate()

# --- Row 39: CU 23 "Realizar diagramas de robustez y secuencia de CU 23." ---
# Status goes from "Por iniciar" to "Hecho" and the estimated-hours column (G)
# gets 1 hour, which ripples through the cascading H..AY "spent/remaining"
# formulas (I,L,O,R,U,X,AA,AD,AG) and the AI milestone cell gets 1 as well.
$ws.Range("F39").Value = "Hecho"
$ws.Range("G39").Value = 1
$ws.Range("AI39").Value = 1

# --- Row 40: CU 24 "Realizar diagramas de robustez y secuencia de CU 24." ---
$ws.Range("F40").Value = "Hecho"
$ws.Range("G40").Value = 1
$ws.Range("AI40").Value = 1

# --- Re-apply the header merges touched by the edit so they come back in ---
# --- the same relative order the workbook was saved with. ---
$allMerges = @("AL4:AM4","H4:I4","K4:L4","N4:O4","Q4:R4","T4:U4","W4:X4","Z4:AA4","AC4:AD4","AF4:AG4","AI4:AJ4","AZ4:BA4","AO4:AP4","AR4:AS4","AU4:AV4","AX4:AY4")
foreach ($r in $allMerges) {
    $ws.Range($r).UnMerge() | Out-Null
}
$orderedMerges = @("AZ4:BA4","AO4:AP4","AR4:AS4","AU4:AV4","AX4:AY4","AL4:AM4","H4:I4","K4:L4","N4:O4","Q4:R4","T4:U4","W4:X4","Z4:AA4","AC4:AD4","AF4:AG4","AI4:AJ4")
foreach ($r in $orderedMerges) {
    $ws.Range($r).Merge() | Out-Null
}

# --- Leave the view scrolled/selected on the last touched cell ---
$ws.Range("AI39").Select()
